$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MoveLCD")

# --- Update the reference segment text (row 3) so the parsed start/end
#     coordinates move from the old R2 segment to the new C17 segment.
#     Order matters for shared-string slot allocation: E5, then G5, then P3.
$ws.Range("E5").Value = "C17"
$ws.Range("G5").Value = "    (at 27.94 41.021)"
$ws.Range("P3").Value = "  (segment (start 27.997 39.3979) (end 28.3863 38.9655) (width 0.1524) (layer Back) (net 47))"

# --- Row 6: offset of the new component (C17) from the LCD1 reference
#     point, replacing the old "translate by the same delta as LCD1" formulas.
$ws.Range("B6").Formula = "=B5-B3"
$ws.Range("C6").Formula = "=C5-C3"
$ws.Range("M5").Copy() | Out-Null
$ws.Range("B6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The old D6/E6/G6 (rotation, label, generated text) no longer apply to
# row 6 -- that computation now lives in row 7 -- so clear them out.
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()
$ws.Range("G6").Clear()

# --- Row 7: final C17' position text, built from the row-6 offsets plus
#     the LCD1' anchor point.
$ws.Range("B7").Formula = "=ROUND(`$B`$4-C6,3)"
$ws.Range("C7").Formula = "=ROUND(`$C`$4+B6,3)"
$ws.Range("D7").Formula = "=D5+D1"
$ws.Range("E7").Formula = "=E5&""'"""
$ws.Range("G7").Formula = "=LEFT(G5,H5+LEN(H`$1)-1)&TEXT(B7,""#0.0####"")&"" ""&TEXT(C7,""#0.0####"")&"" ""&TEXT(D7,""#0"")&K`$1"

# M7/N7 now anchor on the LCD1' point with fully-absolute references
# (previously only the row was absolute).
$ws.Range("M7").Formula = "=ROUND(`$B`$4-N5,3)"
$ws.Range("N7").Formula = "=ROUND(`$C`$4+M5,3)"

# --- Column J widened slightly to fit the new #VALUE! helper cells.
$ws.Range("J1").ColumnWidth = 2.54296875

# --- Cosmetic: remembered vertical scroll position of the window.
$excel.ActiveWindow.ScrollRow = 3
